# Commit: "Suppression 1 alerte + modif alerts"
# The row for "Echappée" on 44836 (25/04→row 25) is removed from the
# planning table; every row below it shifts up by one. We reproduce
# that by deleting worksheet row 25 outright (Excel's Delete Entire Row),
# which removes the cell content and re-indexes all following rows -
# exactly matching the shared-strings usage count drop (76 -> 74) and
# the new dimension (A1:C37 -> A1:C36) seen in the target file.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows("25").Delete()

# Reflect the editor's final on-screen state: after deleting the row,
# the whole new row 25 (now "Coin du Balai" / 44843) is selected and the
# view has scrolled down so row 17 is at the top.
$ws.Range("A25:XFD25").Select()
$excel.ActiveWindow.ScrollRow = 17
